$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.455.25"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "3.605.46"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.17"
$ws.Range("E5").Value = "  -2.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.39"
$ws.Range("E6").Value = "  -3.21%  "

$ws.Range("D7").Value = "3.601.58"
$ws.Range("E7").Value = "  -0.73%  "

$ws.Range("E8").Value = "  -2.46%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  +3.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.660"
$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.91"
$ws.Range("E12").Value = "  -4.92%  "

$ws.Range("E13").Value = "  +7.12%  "

$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("D15").Value = "4.180.12"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.78"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "3.597.61"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "70.328.27"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("E19").Value = "  -0.88%  "

$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.88"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.40"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("E24").Value = "  -9.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.45"
$ws.Range("E25").Value = "  +5.56%  "

$ws.Range("E26").Value = "  -2.58%  "

$ws.Range("E27").Value = "  -6.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  -4.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  -2.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.18"
$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("E31").Value = "  -4.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.18"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.78"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "567.32"
$ws.Range("E35").Value = "  -9.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.77"
$ws.Range("E36").Value = "  -3.98%  "

$ws.Range("D37").Value = "0.0₃0811"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.395"
$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +14.06%  "

$ws.Range("E41").Value = "  +5.17%  "

$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.137"
$ws.Range("E43").Value = "  -6.31%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.01"
$ws.Range("E44").Value = "  -4.61%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.54"
$ws.Range("E45").Value = "  +6.07%  "

$ws.Range("D46").Value = "3.217.19"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.70"
$ws.Range("E48").Value = "  +4.87%  "

$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.21"
$ws.Range("E51").Value = "  -4.47%  "
